$d = $word.ActiveDocument

# The paragraph "<id>p144r_1</id>" is currently split across three runs:
#   <id>    (Courier New, color 7f6000, sz 18)
#   p144r_1 (color 000000)
#   </id>   (Courier New, color 7f6000, sz 18)
# The edit merges them into a single run carrying the first run's
# formatting, with the combined text "<id>p144r_1</id>".
#
# Scope the search to the specific paragraph (index 5, 1-based) that
# contains this text so we only touch the intended run and don't risk
# matching the similarly-shaped "<id>fig_p144r_1</id>" elsewhere.
$p = $d.Paragraphs(5)
$r = $p.Range
$r.Find.ClearFormatting()
$r.Find.Execute("<id>p144r_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p144r_1</id>", 2)
